# Applies the "jesuita-entrada_totals" update:
#  - count/date refresh for the top few aggregate rows (rows 2,3,4,5,6,7,8)
#  - a brand-new "Lisboa, Arroios" row inserted at row 22 (pushes the rest
#    of the table down by one row, 73 -> 74 data rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $val) {
    # Columns C/D store digit-only date codes as literal TEXT (not numbers).
    # A plain .Value assignment of a numeric-looking string gets silently
    # re-typed as a Number by Excel's auto-detection, so force Text format,
    # assign, then drop back to the unstyled "Normal" cell style so no stray
    # formatting is left behind (matches how the rest of the sheet looks).
    $c = $sheet.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Row 2: "?" aggregate count changes -----------------------------------
$ws.Range("B2").Value = 493

# --- Row 3: Coimbra count changes ------------------------------------------
$ws.Range("B3").Value = 56

# --- Row 4: was Lisboa -> becomes Paris ------------------------------------
$ws.Range("A4").Value = "Paris"
$ws.Range("B4").Value = 42
Set-TextCell $ws "C4" "15340815"
Set-TextCell $ws "D4" "17590310"

# --- Row 5: was Paris -> becomes Roma ---------------------------------------
$ws.Range("A5").Value = "Roma"
$ws.Range("B5").Value = 41
Set-TextCell $ws "C5" "15400927"
Set-TextCell $ws "D5" "17560709"

# --- Row 6: was Roma -> becomes Lisboa --------------------------------------
$ws.Range("A6").Value = "Lisboa"
$ws.Range("B6").Value = 38
Set-TextCell $ws "C6" "15460000"
Set-TextCell $ws "D6" "17521123"

# --- Row 7: Goa count changes ------------------------------------------------
$ws.Range("B7").Value = 33

# --- Row 8: Évora count changes -----------------------------------------------
$ws.Range("B8").Value = 26

# --- Insert new row 22 "Lisboa, Arroios" (shifts old rows 22-73 to 23-74) ---
$ws.Rows.Item(22).Insert()

# Give the new row's label cell the same style as its neighbours (bold,
# bordered, centered) by copying formats down from the row above.
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A22").Value = "Lisboa, Arroios"
$ws.Range("B22").Value = 6
Set-TextCell $ws "C22" "17450423"
Set-TextCell $ws "D22" "17530612"
